$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (VP), C (FN), D (FP), E (Precision), F (Recall)
# for rows 2-12.
$data = @{
    2  = @(5, 6, 2, 0.7142857142857143, 0.4545454545454545)
    3  = @(15, 6, 6, 0.7142857142857143, 0.7142857142857143)
    4  = @(15, 6, 6, 0.7142857142857143, 0.7142857142857143)
    5  = @(21, 7, 8, 0.7241379310344828, 0.75)
    6  = @(21, 7, 15, 0.5833333333333334, 0.75)
    7  = @(24, 9, 15, 0.6153846153846154, 0.7272727272727273)
    8  = @(24, 11, 15, 0.6153846153846154, 0.6857142857142857)
    9  = @(26, 12, 18, 0.5909090909090909, 0.6842105263157895)
    10 = @(26, 12, 18, 0.5909090909090909, 0.6842105263157895)
    11 = @(26, 12, 32, 0.4482758620689655, 0.6842105263157895)
    12 = @(26, 12, 32, 0.4482758620689655, 0.6842105263157895)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
}
